$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cells in this sheet store every value (coin name, link, price, volume)
# as plain text, even when the text looks like a number (e.g. "0.999" or
# "66.123.98"). Excel normally auto-converts such text to a real number when
# you assign it through .Value, so we first force a text number format on
# every cell we are about to touch (kept as separate, non-comma Range calls
# because applying NumberFormat through a single multi-area/union Range is
# unreliable), write the values, and finally reset the style back to the
# workbook default so no stray formatting is left behind.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("B32").NumberFormat = "@"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("C48").NumberFormat = "@"

$ws.Range('D2').Value = '66.123.98'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '3.566.59'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('D5').Value = '607.07'
$ws.Range('D6').Value = '145.32'
$ws.Range('E6').Value = '  +1.75%  '
$ws.Range('D7').Value = '3.567.30'
$ws.Range('E7').Value = '  +2.64%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  +3.53%  '
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  -3.40%  '
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '4.170.83'
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '0.0000208'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '30.03'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').Value = '3.556.71'
$ws.Range('E16').Value = '  +2.34%  '
$ws.Range('D17').Value = '66.211.89'
$ws.Range('E17').Value = '  +0.14%  '
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').Value = '11.41'
$ws.Range('E19').Value = '  +9.56%  '
$ws.Range('D21').Value = '14.88'
$ws.Range('E21').Value = '  +1.27%  '
$ws.Range('D22').Value = '429.79'
$ws.Range('E22').Value = '  +2.42%  '
$ws.Range('D23').Value = '0.615'
$ws.Range('E23').Value = '  +4.58%  '
$ws.Range('D24').Value = '79.16'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('D25').Value = '3.708.08'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('E27').Value = '  +3.89%  '
$ws.Range('D28').Value = '2.51'
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('D29').Value = '7.96'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').Value = '9.10'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.47'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '25.60'
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').Value = '3.559.11'
$ws.Range('E34').Value = '  +2.50%  '
$ws.Range('E35').Value = '  -5.91%  '
$ws.Range('E37').Value = '  +2.34%  '
$ws.Range('D38').Value = '7.87'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('D39').Value = '5.62'
$ws.Range('E39').Value = '  +1.43%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').Value = '175.09'
$ws.Range('E41').Value = '  +3.34%  '
$ws.Range('D42').Value = '0.0849'
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').Value = '5.22'
$ws.Range('E43').Value = '  +2.48%  '
$ws.Range('D44').Value = '0.897'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').Value = '46.15'
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '25.79'
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('D49').Value = '2.40'
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D50').Value = '23.52'
$ws.Range('E50').Value = '  +9.59%  '
$ws.Range('D51').Value = '7.14'
$ws.Range('E51').Value = '  +0.33%  '

$ws.Range("D2:E51").Style = "Normal"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Style = "Normal"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Style = "Normal"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Style = "Normal"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Style = "Normal"

Write-Host "Applied cryptos list update"
